$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").ClearContents()
$ws.Range("A2").Value = "Maryland 2021 Monthly Dispensary Edibles Sales"
Write-Host "done"
